# Insert a new daily price record for "Feria Lagunitas de Puerto Montt - Perejil"
# right before the existing row 376, pushing all subsequent rows down by one
# (old row 376 becomes 377, ..., old row 413 becomes 414).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 376..413 down to 377..414, leaving a blank row 376 to fill in.
$ws.Rows.Item(376).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(376, 1).Value  = 4
$ws.Cells.Item(376, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(376, 3).Value  = 'Los Lagos'
$ws.Cells.Item(376, 4).Value  = 45106
$ws.Cells.Item(376, 5).Value  = 10
$ws.Cells.Item(376, 6).Value  = 100112044
$ws.Cells.Item(376, 7).Value  = 'Perejil'
$ws.Cells.Item(376, 8).Value  = 'Sin especificar'
$ws.Cells.Item(376, 9).Value  = 'Primera'
$ws.Cells.Item(376, 10).Value = 60
$ws.Cells.Item(376, 11).Value = 7000
$ws.Cells.Item(376, 12).Value = 7000
$ws.Cells.Item(376, 13).Value = 7000
$ws.Cells.Item(376, 14).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(376, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(376, 16).Value = 2333
$ws.Cells.Item(376, 17).Value = 3
$ws.Cells.Item(376, 18).Value = 'Hortaliza'
